$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B14").Value = "f"
$ws.Range("C14").Value = 27

$ws.Range("B15").Value = "f"
$ws.Range("C15").Value = 29

$ws.Range("B16").Value = "f"
$ws.Range("C16").Value = 27

$ws.Range("C17").Select()
